# The source workbook tracks daily "Cebollín" price records for Feria
# Lagunitas de Puerto Montt, one row per day. This commit adds a missing
# daily record, inserted right before the existing 2021-10-26 (serial
# 44330) row so the data stays in (roughly) chronological order — which
# pushes every subsequent record down by one row (old row 76 -> new row
# 77, ..., old row 171 -> new row 172).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 76; Excel shifts rows 76:171 down to
# 77:172 and the sheet's used range/dimension grows to A1:R172
# automatically.
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A76").Value = 4
$ws.Range("B76").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C76").Value = "Los Lagos"
$ws.Range("D76").Value = 44483
$ws.Range("E76").Value = 10
$ws.Range("F76").Value = 100112037
$ws.Range("G76").Value = "Cebollín"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 90
$ws.Range("K76").Value = 5500
$ws.Range("L76").Value = 5500
$ws.Range("M76").Value = 5500
$ws.Range("N76").Value = "$/paquete 36 unidades"
$ws.Range("O76").Value = "Región Metropolitana"
$ws.Range("P76").Value = 153
$ws.Range("Q76").Value = 36
$ws.Range("R76").Value = "Hortaliza"
